$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for Price/Volume columns so numeric-looking values
# (e.g. "0.999") are preserved as text rather than being converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "45.282.14"

$ws.Range("D3").Value = "2.421.28"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "319.34"
$ws.Range("E5").Value = "  +3.69%  "

$ws.Range("D6").Value = "103.05"
$ws.Range("E6").Value = "  +2.14%  "

$ws.Range("D7").Value = "0.514"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +5.11%  "

$ws.Range("D10").Value = "35.48"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("D13").Value = "18.24"
$ws.Range("E13").Value = "  -3.51%  "

$ws.Range("E14").Value = "  +1.29%  "

$ws.Range("D15").Value = "2.799.99"
$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").Value = "2.417.71"
$ws.Range("E16").Value = "  -2.09%  "

$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "45.173.43"
$ws.Range("E18").Value = "  +2.33%  "

$ws.Range("D19").Value = "12.26"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "6.33"
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").Value = "0.0₃0921"
$ws.Range("E21").Value = "  +1.64%  "

$ws.Range("D22").Value = "69.88"
$ws.Range("E22").Value = "  +1.86%  "

$ws.Range("D23").Value = "244.39"
$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("E24").Value = "  -2.06%  "

$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +1.75%  "

$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -6.18%  "

$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  +0.77%  "

$ws.Range("D30").Value = "49.47"
$ws.Range("E30").Value = "  +2.35%  "

$ws.Range("D31").Value = "32.93"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").Value = "20.20"
$ws.Range("E32").Value = "  +8.12%  "

$ws.Range("D33").Value = "0.126"
$ws.Range("E33").Value = "  +7.12%  "

$ws.Range("D34").Value = "5.23"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("D36").Value = "0.0761"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("D38").Value = "4.44"
$ws.Range("E38").Value = "  -0.73%  "

$ws.Range("D39").Value = "128.40"
$ws.Range("E39").Value = "  -2.54%  "

$ws.Range("E40").Value = "  -0.70%  "

$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("D43").Value = "20.48"
$ws.Range("E43").Value = "  -5.24%  "

$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("D45").Value = "1.942.74"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("E46").Value = "  -2.65%  "

$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("E48").Value = "  +8.44%  "

$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("D50").Value = "76.90"
$ws.Range("E50").Value = "  +4.25%  "

$ws.Range("E51").Value = "  +5.10%  "

# Restore default cell style (the text number-format above is only a transient
# aid for entry; the original cells had no explicit style override).
$ws.Range("D2:E51").Style = "Normal"
